# Apply updated cryptocurrency price/volume figures (columns D and E) for rows 2-51.
# Column D cells that are purely numeric-looking strings get a leading apostrophe
# so Excel stores them as text (matching the original inlineStr cells) instead of
# auto-converting them to a Number type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.863.29'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '2.029.52'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''227.34'
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("D7").Value = '''59.56'
$ws.Range("E7").Value = '  +2.65%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = '''0.384'
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("E11").Value = '  +0.76%  '
$ws.Range("D12").Value = '''14.57'
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").Value = '2.330.70'
$ws.Range("E13").Value = '  -1.02%  '
$ws.Range("D14").Value = '''21.14'
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("D15").Value = '''0.760'
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("E16").Value = '  -1.69%  '
$ws.Range("D17").Value = '2.028.93'
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("D18").Value = '37.768.19'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("E19").Value = '  -1.75%  '
$ws.Range("D20").Value = '''70.03'
$ws.Range("E20").Value = '  +0.67%  '
$ws.Range("E21").Value = '  -0.97%  '
$ws.Range("D22").Value = '''224.88'
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '''2.39'
$ws.Range("E24").Value = '  -2.27%  '
$ws.Range("D25").Value = '''2.19'
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("D26").Value = '''9.23'
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").Value = '''165.31'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("E28").Value = '  -2.73%  '
$ws.Range("D29").Value = '''18.94'
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("E30").Value = '  -4.12%  '
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("D33").Value = '''2.09'
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("E35").Value = '  -1.27%  '
$ws.Range("E36").Value = '  +6.92%  '
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("E38").Value = '  -2.04%  '
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("D40").Value = '1.518.43'
$ws.Range("E40").Value = '  +2.45%  '
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").Value = '''96.63'
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("D43").Value = '''16.78'
$ws.Range("E43").Value = '  +1.31%  '
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").Value = '''0.0918'
$ws.Range("E45").Value = '  -1.65%  '
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("D47").Value = '''4.05'
$ws.Range("E47").Value = '  -3.72%  '
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("D50").Value = '''7.03'
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("D51").Value = '2.218.42'
$ws.Range("E51").Value = '  -1.10%  '
